# Auto-generated Excel COM-interop edit script
# Applies updated crypto price / volume(1h) values, and fixes the swapped
# EnergySwap / Algorand rows (47 and 48), per the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell text updates --------------------------------------------------
# All updated cells (B/C "Coin"/"Link" swap for rows 47-48, plus the
# D "Price" / E "Volume(1h)" columns) are stored as text in the workbook,
# even though many of the D/E values look like numbers (e.g. "263.04",
# "1.001", "  +0.07%  "). Setting .Value directly on a numeric-looking
# string would make Excel auto-convert the cell to a real number, which
# would change the cell's underlying type and not match the source data.
# To keep every cell as text (matching the original file), we temporarily
# force the cell's number format to Text ("@") before assigning the value,
# then restore the cell's original style afterwards so no stray formatting
# is left behind. This is harmless for the plain (non-numeric-looking)
# B/C strings too.
$updates = @(
    @{Cell="D2"; Value="26.621.61"},
    @{Cell="E2"; Value="  +0.07%  "},
    @{Cell="D3"; Value="1.851.55"},
    @{Cell="E3"; Value="  -0.16%  "},
    @{Cell="D4"; Value="1.001"},
    @{Cell="E4"; Value="  +0.08%  "},
    @{Cell="D5"; Value="263.04"},
    @{Cell="E5"; Value="  -0.56%  "},
    @{Cell="E6"; Value="  +0.09%  "},
    @{Cell="D7"; Value="0.5336"},
    @{Cell="E7"; Value="  +2.03%  "},
    @{Cell="D8"; Value="0.3163"},
    @{Cell="E8"; Value="  -3.60%  "},
    @{Cell="D9"; Value="0.06958"},
    @{Cell="E9"; Value="  +2.31%  "},
    @{Cell="D10"; Value="18.88"},
    @{Cell="E10"; Value="  +0.25%  "},
    @{Cell="D11"; Value="0.7711"},
    @{Cell="E11"; Value="  -0.79%  "},
    @{Cell="D12"; Value="0.07836"},
    @{Cell="E12"; Value="  +0.92%  "},
    @{Cell="D13"; Value="1.851.17"},
    @{Cell="E13"; Value="  -0.29%  "},
    @{Cell="D14"; Value="89.81"},
    @{Cell="D15"; Value="5.054"},
    @{Cell="E15"; Value="  +0.54%  "},
    @{Cell="D16"; Value="14.14"},
    @{Cell="D17"; Value="1.001"},
    @{Cell="E17"; Value="  +0.03%  "},
    @{Cell="D18"; Value="0.000007971"},
    @{Cell="E18"; Value="  -0.06%  "},
    @{Cell="E19"; Value="  +0.04%  "},
    @{Cell="D20"; Value="26.643.80"},
    @{Cell="E20"; Value="  +0.04%  "},
    @{Cell="D21"; Value="2.095.41"},
    @{Cell="E21"; Value="  +0.50%  "},
    @{Cell="D22"; Value="4.654"},
    @{Cell="E22"; Value="  +0.28%  "},
    @{Cell="D23"; Value="6.031"},
    @{Cell="E23"; Value="  +0.52%  "},
    @{Cell="D24"; Value="9.362"},
    @{Cell="E24"; Value="  -2.05%  "},
    @{Cell="D25"; Value="2.214"},
    @{Cell="E25"; Value="  +0.74%  "},
    @{Cell="D26"; Value="141.76"},
    @{Cell="E26"; Value="  -1.83%  "},
    @{Cell="D27"; Value="1.692"},
    @{Cell="E27"; Value="  +0.88%  "},
    @{Cell="D28"; Value="17.13"},
    @{Cell="E28"; Value="  +0.68%  "},
    @{Cell="D29"; Value="111.82"},
    @{Cell="E29"; Value="  -0.48%  "},
    @{Cell="D30"; Value="4.309"},
    @{Cell="E30"; Value="  +2.54%  "},
    @{Cell="D31"; Value="0.08782"},
    @{Cell="E31"; Value="  +0.28%  "},
    @{Cell="D32"; Value="4.121"},
    @{Cell="E32"; Value="  -0.88%  "},
    @{Cell="D33"; Value="0.04859"},
    @{Cell="E33"; Value="  +0.55%  "},
    @{Cell="D34"; Value="0.7387"},
    @{Cell="E34"; Value="  +3.39%  "},
    @{Cell="D35"; Value="1.140"},
    @{Cell="E35"; Value="  +0.10%  "},
    @{Cell="D36"; Value="2.889"},
    @{Cell="E36"; Value="  +0.95%  "},
    @{Cell="E37"; Value="  +0.00%  "},
    @{Cell="D38"; Value="2.349"},
    @{Cell="E38"; Value="  +6.49%  "},
    @{Cell="D39"; Value="0.01737"},
    @{Cell="E39"; Value="  -2.77%  "},
    @{Cell="D40"; Value="0.4836"},
    @{Cell="E40"; Value="  -0.96%  "},
    @{Cell="D41"; Value="0.9076"},
    @{Cell="E41"; Value="  +0.59%  "},
    @{Cell="D42"; Value="108.87"},
    @{Cell="E42"; Value="  -3.59%  "},
    @{Cell="D43"; Value="5.921"},
    @{Cell="E43"; Value="  -2.78%  "},
    @{Cell="E44"; Value="  +0.09%  "},
    @{Cell="D45"; Value="7.698"},
    @{Cell="E45"; Value="  -0.40%  "},
    @{Cell="D46"; Value="0.4203"},
    @{Cell="E46"; Value="  +0.14%  "},
    @{Cell="B47"; Value="Algorand"},
    @{Cell="C47"; Value="https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"},
    @{Cell="D47"; Value="0.1252"},
    @{Cell="E47"; Value="  +0.56%  "},
    @{Cell="B48"; Value="EnergySwap"},
    @{Cell="C48"; Value="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"},
    @{Cell="D48"; Value="9.064"},
    @{Cell="E48"; Value="  -0.94%  "},
    @{Cell="D49"; Value="35.05"},
    @{Cell="E49"; Value="  +0.11%  "},
    @{Cell="D51"; Value="0.8978"},
    @{Cell="E51"; Value="  +1.40%  "}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = $origStyle
}
